$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest scrape.
# Some Price values look numeric to Excel (e.g. "674.03"); force those
# cells to remain plain text, matching the original inline-string data,
# by setting the number format to Text ("@") before assigning the value.

$ws.Range("D2").Value = "69.174.96"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "3.669.26"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "674.03"
$ws.Range("E5").Value = "  -1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.35"
$ws.Range("E6").Value = "  -3.24%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -1.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.145"
$ws.Range("E9").Value = "  -1.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.91"
$ws.Range("E10").Value = "  -5.89%  "

$ws.Range("E11").Value = "  -2.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000230"
$ws.Range("E12").Value = "  -3.64%  "

$ws.Range("D13").Value = "4.288.35"
$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.10"
$ws.Range("E14").Value = "  -4.28%  "

$ws.Range("D15").Value = "3.656.57"
$ws.Range("E15").Value = "  -0.95%  "

$ws.Range("D16").Value = "69.170.50"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.96"
$ws.Range("E18").Value = "  -1.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  -3.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.26"
$ws.Range("E20").Value = "  -3.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.92"
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.645"
$ws.Range("E22").Value = "  -3.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.78"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("D24").Value = "3.816.38"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.84"
$ws.Range("E26").Value = "  -5.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  -8.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.94"
$ws.Range("E28").Value = "  -6.79%  "

$ws.Range("E29").Value = "  -2.50%  "

$ws.Range("E30").Value = "  -6.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.57"
$ws.Range("E31").Value = "  -4.26%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.80"
$ws.Range("E33").Value = "  -1.15%  "

$ws.Range("E34").Value = "  -5.81%  "

$ws.Range("D35").Value = "3.660.04"
$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("E36").Value = "  -4.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.09"
$ws.Range("E37").Value = "  -5.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.12"
$ws.Range("E38").Value = "  -3.82%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.19"
$ws.Range("E41").Value = "  -2.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0892"
$ws.Range("E42").Value = "  -4.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "171.25"
$ws.Range("E43").Value = "  +7.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.939"
$ws.Range("E44").Value = "  -1.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.46"
$ws.Range("E45").Value = "  -1.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000274"
$ws.Range("E46").Value = "  -5.17%  "

$ws.Range("E47").Value = "  -6.82%  "

$ws.Range("E50").Value = "  -4.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.73"
$ws.Range("E51").Value = "  -3.84%  "

# Row 48/49: ONDO and InjectiveProtocol swapped position in the ranking,
# each also picking up refreshed Price/Volume figures.
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.19"
$ws.Range("E48").Value = "  -8.98%  "

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.27"
$ws.Range("E49").Value = "  -6.74%  "
